# Auto-generated edit: update currentAveragePrice / Leve price & profit
# columns (H-N) for the affected Leve rows across sheets, per the scheduled
# market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2926.08
$ws.Range("I98").Value = 3181.5
$ws.Range("J98").Value = 2269.2856
$ws.Range("K98").Value = 3181.5
$ws.Range("L98").Value = 2269.2856
$ws.Range("M98").Value = -1683.5
$ws.Range("N98").Value = -5265.2856
$ws.Range("H100").Value = 4274.375
$ws.Range("I100").Value = 3425
$ws.Range("J100").Value = 4784
$ws.Range("K100").Value = 3425
$ws.Range("L100").Value = 4784
$ws.Range("M100").Value = -2884
$ws.Range("N100").Value = -5866
$ws.Range("H106").Value = 2440.5
$ws.Range("I106").Value = 800.8333
$ws.Range("K106").Value = 800.8333
$ws.Range("M106").Value = -169.8333
$ws.Range("H107").Value = 551.0476
$ws.Range("I107").Value = 444.66666
$ws.Range("J107").Value = 817
$ws.Range("K107").Value = 444.66666
$ws.Range("L107").Value = 817
$ws.Range("M107").Value = 1475.33334
$ws.Range("N107").Value = -4657
$ws.Range("H113").Value = 4159.9287
$ws.Range("I113").Value = 4025.3635
$ws.Range("J113").Value = 4653.3335
$ws.Range("K113").Value = 4025.3635
$ws.Range("L113").Value = 4653.3335
$ws.Range("M113").Value = -771.3634999999999
$ws.Range("N113").Value = -11161.3335
$ws.Range("H115").Value = 1293.3334
$ws.Range("I115").Value = 510
$ws.Range("J115").Value = 2272.5
$ws.Range("K115").Value = 1530
$ws.Range("L115").Value = 6817.5
$ws.Range("M115").Value = 37
$ws.Range("N115").Value = -9951.5
$ws.Range("H118").Value = 2350.2942
$ws.Range("I118").Value = 328.8889
$ws.Range("J118").Value = 4624.375
$ws.Range("K118").Value = 986.6667
$ws.Range("L118").Value = 13873.125
$ws.Range("M118").Value = 670.3333
$ws.Range("N118").Value = -17187.125
$ws.Range("H122").Value = 2926.08
$ws.Range("I122").Value = 3181.5
$ws.Range("J122").Value = 2269.2856
$ws.Range("K122").Value = 9544.5
$ws.Range("L122").Value = 6807.8568
$ws.Range("M122").Value = -7094.5
$ws.Range("N122").Value = -11707.8568
$ws.Range("H132").Value = 4547880
$ws.Range("I132").Value = 5130366
$ws.Range("J132").Value = 4488
$ws.Range("K132").Value = 15391098
$ws.Range("L132").Value = 13464
$ws.Range("M132").Value = -15388568
$ws.Range("N132").Value = -18524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6738.739
$ws.Range("I32").Value = 5818.9375
$ws.Range("K32").Value = 5818.9375
$ws.Range("M32").Value = -5531.9375
$ws.Range("H97").Value = 1455.9333
$ws.Range("I97").Value = 1347.0714
$ws.Range("K97").Value = 1347.0714
$ws.Range("M97").Value = -851.0714
$ws.Range("H110").Value = 1696.8096
$ws.Range("I110").Value = 701.6667
$ws.Range("J110").Value = 3023.6667
$ws.Range("K110").Value = 701.6667
$ws.Range("L110").Value = 3023.6667
$ws.Range("M110").Value = 1343.3333
$ws.Range("N110").Value = -7113.6667
$ws.Range("H132").Value = 1759.5471
$ws.Range("I132").Value = 1309.6522
$ws.Range("K132").Value = 3928.9566
$ws.Range("M132").Value = -1398.9566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 65527
$ws.Range("J38").Value = 65527
$ws.Range("L38").Value = 65527
$ws.Range("N38").Value = -66359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1524.2858
$ws.Range("I16").Value = 941.4167
$ws.Range("J16").Value = 2301.4443
$ws.Range("K16").Value = 941.4167
$ws.Range("L16").Value = 2301.4443
$ws.Range("M16").Value = -654.4167
$ws.Range("N16").Value = -2875.4443
$ws.Range("H62").Value = 4394.364
$ws.Range("I62").Value = 3416.875
$ws.Range("J62").Value = 7001
$ws.Range("K62").Value = 3416.875
$ws.Range("L62").Value = 7001
$ws.Range("M62").Value = -2792.875
$ws.Range("N62").Value = -8249
$ws.Range("H65").Value = 4394.364
$ws.Range("I65").Value = 3416.875
$ws.Range("J65").Value = 7001
$ws.Range("K65").Value = 17084.375
$ws.Range("L65").Value = 35005
$ws.Range("M65").Value = -13964.375
$ws.Range("N65").Value = -41245
$ws.Range("H86").Value = 8455.333000000001
$ws.Range("I86").Value = 7399.769
$ws.Range("J86").Value = 11199.8
$ws.Range("K86").Value = 7399.769
$ws.Range("L86").Value = 11199.8
$ws.Range("M86").Value = -6276.769
$ws.Range("N86").Value = -13445.8
$ws.Range("H89").Value = 8455.333000000001
$ws.Range("I89").Value = 7399.769
$ws.Range("J89").Value = 11199.8
$ws.Range("K89").Value = 36998.845
$ws.Range("L89").Value = 55999
$ws.Range("M89").Value = -31382.845
$ws.Range("N89").Value = -67231
$ws.Range("H99").Value = 2514.1428
$ws.Range("I99").Value = 1025
$ws.Range("J99").Value = 4499.6665
$ws.Range("K99").Value = 1025
$ws.Range("L99").Value = 4499.6665
$ws.Range("M99").Value = 473
$ws.Range("N99").Value = -7495.6665
$ws.Range("H107").Value = 1433.6
$ws.Range("I107").Value = 514.46155
$ws.Range("J107").Value = 3140.5715
$ws.Range("K107").Value = 514.46155
$ws.Range("L107").Value = 3140.5715
$ws.Range("M107").Value = 1405.53845
$ws.Range("N107").Value = -6980.5715
$ws.Range("H113").Value = 1524.2858
$ws.Range("I113").Value = 941.4167
$ws.Range("J113").Value = 2301.4443
$ws.Range("K113").Value = 941.4167
$ws.Range("L113").Value = 2301.4443
$ws.Range("M113").Value = 1228.5833
$ws.Range("N113").Value = -6641.4443
$ws.Range("H122").Value = 4364.727
$ws.Range("I122").Value = 3251.5
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 9754.5
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -7304.5
$ws.Range("N122").Value = -26900.0005
$ws.Range("H126").Value = 2514.1428
$ws.Range("I126").Value = 1025
$ws.Range("J126").Value = 4499.6665
$ws.Range("K126").Value = 3075
$ws.Range("L126").Value = 13498.9995
$ws.Range("M126").Value = -605
$ws.Range("N126").Value = -18438.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2133.6667
$ws.Range("J31").Value = 3334
$ws.Range("L31").Value = 10002
$ws.Range("N31").Value = -10578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4349324
$ws.Range("I7").Value = 6250933.5
$ws.Range("J7").Value = 2787.8572
$ws.Range("K7").Value = 6250933.5
$ws.Range("L7").Value = 2787.8572
$ws.Range("M7").Value = -6250821.5
$ws.Range("N7").Value = -3011.8572
$ws.Range("H16").Value = 52632384
$ws.Range("I16").Value = 55556376
$ws.Range("K16").Value = 55556376
$ws.Range("M16").Value = -55556206
$ws.Range("H122").Value = 2297.4666
$ws.Range("I122").Value = 1872.15
$ws.Range("K122").Value = 5616.450000000001
$ws.Range("M122").Value = -3166.450000000001
$ws.Range("H126").Value = 4349324
$ws.Range("I126").Value = 6250933.5
$ws.Range("J126").Value = 2787.8572
$ws.Range("K126").Value = 18752800.5
$ws.Range("L126").Value = 8363.571599999999
$ws.Range("M126").Value = -18750330.5
$ws.Range("N126").Value = -13303.5716
$ws.Range("H132").Value = 4617.386
$ws.Range("I132").Value = 1516.0344
$ws.Range("K132").Value = 4548.1032
$ws.Range("M132").Value = -2018.1032

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 29250
$ws.Range("J46").Value = 29250
$ws.Range("L46").Value = 29250
$ws.Range("N46").Value = -29712
$ws.Range("H107").Value = 1292.2858
$ws.Range("I107").Value = 599.1111
$ws.Range("J107").Value = 2540
$ws.Range("K107").Value = 1797.3333
$ws.Range("L107").Value = 7620
$ws.Range("M107").Value = 122.6667000000002
$ws.Range("N107").Value = -11460
$ws.Range("H126").Value = 2633321.8
$ws.Range("J126").Value = 9093800
$ws.Range("L126").Value = 27281400
$ws.Range("N126").Value = -27286340
$ws.Range("H132").Value = 120446.33
$ws.Range("I132").Value = 145734.33
$ws.Range("J132").Value = 28610.947
$ws.Range("K132").Value = 437202.99
$ws.Range("L132").Value = 85832.841
$ws.Range("M132").Value = -434672.99
$ws.Range("N132").Value = -90892.841
$ws.Range("H134").Value = 29250
$ws.Range("J134").Value = 29250
$ws.Range("L134").Value = 87750
$ws.Range("N134").Value = -92820
